$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.009796047210693359
$ws.Range("C2").Value = 0.02279024124145508
$ws.Range("D2").Value = 0.00312647819519043
$ws.Range("E2").Value = 0.01762537956237793
$ws.Range("F2").Value = 0.003726291656494141
$ws.Range("G2").Value = 0.05270538330078125
$ws.Range("H2").Value = 0.009382867813110351
$ws.Range("I2").Value = 0.02599706649780274
$ws.Range("J2").Value = 0.00970306396484375
$ws.Range("K2").Value = 0.01475033760070801
$ws.Range("L2").Value = 0.00700526237487793
$ws.Range("M2").Value = 0.0132171630859375
$ws.Range("B3").Value = 0.05377283096313477
$ws.Range("C3").Value = 0.02862896919250488
$ws.Range("D3").Value = 0.008992338180541992
$ws.Range("E3").Value = 0.01143684387207031
$ws.Range("F3").Value = 0.0157198429107666
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0.06583108901977539
$ws.Range("I3").Value = 0.02456555366516113
$ws.Range("J3").Value = 0.03796019554138184
$ws.Range("K3").Value = 0.01563382148742676
$ws.Range("L3").Value = 0.01526918411254883
$ws.Range("M3").Value = 0.01136612892150879
$ws.Range("B4").Value = 0.02357845306396484
$ws.Range("C4").Value = 0.01424446105957031
$ws.Range("D4").Value = 0.007114553451538086
$ws.Range("E4").Value = 0.007086801528930664
$ws.Range("F4").Value = 0.02898449897766113
$ws.Range("G4").Value = 0.007525300979614258
$ws.Range("H4").Value = 0.0185636043548584
$ws.Range("I4").Value = 0.01203804016113281
$ws.Range("J4").Value = 0.01062946319580078
$ws.Range("K4").Value = 0.01005759239196777
$ws.Range("L4").Value = 0.03003277778625488
$ws.Range("M4").Value = 0.007064485549926757
$ws.Range("B5").Value = 0.01185088157653809
$ws.Range("C5").Value = 0.01421489715576172
$ws.Range("D5").Value = 0.01569652557373047
$ws.Range("E5").Value = 0
$ws.Range("H5").Value = 0.01393203735351562
$ws.Range("I5").Value = 0.01024899482727051
$ws.Range("J5").Value = 0.01558284759521484
$ws.Range("K5").Value = 0.001567840576171875
$ws.Range("B6").Value = 0.2626073837280273
$ws.Range("C6").Value = 0.04109110832214356
$ws.Range("D6").Value = 0.1431647300720215
$ws.Range("E6").Value = 0.03231921195983887
$ws.Range("F6").Value = 0.1126039505004883
$ws.Range("G6").Value = 0.02222952842712402
$ws.Range("H6").Value = 0.3756266117095947
$ws.Range("I6").Value = 0.05570282936096192
$ws.Range("J6").Value = 0.2094531536102295
$ws.Range("K6").Value = 0.03325514793395996
$ws.Range("L6").Value = 0.1269631385803223
$ws.Range("M6").Value = 0.02043871879577637
